$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update the Runmode for the "Forgot Password" row (row 4, column C)
# from "NO" to a new distinct string "YES" (uppercase, separate from existing "Yes").
$ws.Range("C4").Value = "YES"

# Move the active selection on the sheet to D8 (as recorded in the sheet view).
$ws.Range("D8").Select()
